$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 33, shifting existing rows 33:38 down to 34:39
$ws.Rows("33:33").Insert()

# Populate the newly inserted row 33 with the new data record
$ws.Cells.Item(33, 1).Value = 2
$ws.Cells.Item(33, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(33, 3).Value = "Coquimbo"
$ws.Cells.Item(33, 4).Value = 44476
$ws.Cells.Item(33, 5).Value = 4
$ws.Cells.Item(33, 6).Value = 100112022
$ws.Cells.Item(33, 7).Value = "Arveja Verde"
$ws.Cells.Item(33, 8).Value = "Perfection"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 500
$ws.Cells.Item(33, 11).Value = 23000
$ws.Cells.Item(33, 12).Value = 24000
$ws.Cells.Item(33, 13).Value = 23500
$ws.Cells.Item(33, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(33, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(33, 16).Value = 940
$ws.Cells.Item(33, 17).Value = 25
$ws.Cells.Item(33, 18).Value = "Hortaliza"
